# Update dashboards - 2026-01-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 - T5YIFR (5yr, 5yr Forward)
$ws.Range("N29").Value = 46049
$ws.Range("Q29").Value = 2.21
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = 2.18
$ws.Range("T29").Value = 2.2
$ws.Range("U29").Value = 2.26

# Row 30 - T10YIE (10yr TIPS)
$ws.Range("N30").Value = 46049
$ws.Range("Q30").Value = 2.34
$ws.Range("R30").Value = 2.32
$ws.Range("S30").Value = 2.32
$ws.Range("T30").Value = 2.31
$ws.Range("U30").Value = 2.34

# Row 47 - DFF (FFR)
$ws.Range("N47").Value = 46048

# Row 48 - DGS2 (2y UST)
$ws.Range("N48").Value = 46048
$ws.Range("Q48").Value = 3.56
$ws.Range("R48").Value = 3.6
$ws.Range("S48").Value = 3.61
$ws.Range("T48").Value = 3.6
$ws.Range("U48").Value = 3.6

# Row 49 - DGS5 (5y UST)
$ws.Range("N49").Value = 46048
$ws.Range("Q49").Value = 3.82
$ws.Range("R49").Value = 3.84
$ws.Range("S49").Value = 3.85
$ws.Range("T49").Value = 3.83
$ws.Range("U49").Value = 3.86

# Row 50 - DGS10 (10y UST)
$ws.Range("N50").Value = 46048
$ws.Range("Q50").Value = 4.22
$ws.Range("R50").Value = 4.24
$ws.Range("S50").Value = 4.26
$ws.Range("T50").Value = 4.26
$ws.Range("U50").Value = 4.3

# Row 52 - DBAA (BAA)
$ws.Range("N52").Value = 46048
$ws.Range("Q52").Value = 5.83
$ws.Range("R52").Value = 5.85
$ws.Range("S52").Value = 5.85
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.95
